$wb = $excel.ActiveWorkbook

# --- Sheet: "Range Status" ---
# Zero out the species counts and clear the percentage column entirely
# for rows 2-7 (Historical, Very Restricted, Restricted, Moderate, Large,
# Very Large).
$wsRange = $wb.Worksheets.Item("Range Status")
for ($r = 2; $r -le 7; $r++) {
    $wsRange.Cells.Item($r, 2).Value = 0
    $wsRange.Cells.Item($r, 3).ClearContents()
}

# --- Sheet: "Species qualification" ---
# "Range Analysis" row (row 5) species count goes to 0.
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Cells.Item(5, 2).Value = 0

# --- Sheet: "High Priority break-up" ---
# Re-summarised breakdown: now only two categories remain (Trend New, IUCN),
# with updated counts/percentages, and the old "Trend Different" / "Range"
# rows are gone.
$wsBreak = $wb.Worksheets.Item("High Priority break-up")

$wsBreak.Cells.Item(2, 2).Value = 4
$wsBreak.Cells.Item(2, 3).Value = 16.7
$wsBreak.Cells.Item(2, 4).Value = 4
$wsBreak.Cells.Item(2, 5).Value = 16.7

$wsBreak.Cells.Item(3, 1).Value = "IUCN"
$wsBreak.Cells.Item(3, 2).Value = 20
$wsBreak.Cells.Item(3, 3).Value = 83.3
$wsBreak.Cells.Item(3, 4).Value = 20
$wsBreak.Cells.Item(3, 5).Value = 83.3

# Remove the old rows 4 and 5 (former "Range" and "IUCN" rows) entirely.
$wsBreak.Range("A4:E5").Delete()
